$wb = $excel.ActiveWorkbook

# --- Metadata sheet ---
$meta = $wb.Worksheets.Item("Metadata")

# Remove the duplicated "Contact" row (old row 11); this shifts everything below up by one.
$meta.Rows("11").Delete()

# Version bump
$meta.Range("B3").Value = "6.0.0"

# Date refresh
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher now has a value
$meta.Range("B9").Value = "Alvearie Team"

# Former duplicate "Contact" row becomes the new "Jurisdiction" row
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

# --- Elements sheet ---
$elements = $wb.Worksheets.Item("Elements")

# Root Extension row gets a real Short/Definition instead of the generic placeholder text
$elements.Range("K2").Value = "Communication Product"
$elements.Range("L2").Value = "Product generating the communication"
